$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    53  = "BB"
    55  = "BB"
    56  = "A"
    57  = "A"
    130 = "BB"
    156 = "BBB"
    164 = "BB"
    166 = "A"
    172 = "A"
    173 = "A"
    174 = "A"
    183 = "BB"
    213 = "A"
    222 = "A"
    241 = "A"
    242 = "A"
    313 = "BB"
    322 = "BB"
    323 = "BB"
    324 = "BB"
    331 = "BB"
    375 = "A"
    380 = "BBB"
}

foreach ($row in $changes.Keys) {
    $ws.Range("B$row").Value = $changes[$row]
}
